$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.558.32'
$ws.Range("E2").Value = '  +1.13%  '

$ws.Range("D3").Value = '1.572.91'
$ws.Range("E3").Value = '  -1.10%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.93'
$ws.Range("E5").Value = '  -0.36%  '

$ws.Range("E6").Value = '  -0.41%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '45.61'
$ws.Range("E8").Value = '  +3.97%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.18'
$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("E10").Value = '  -1.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0888'
$ws.Range("E12").Value = '  -0.14%  '

$ws.Range("D13").Value = '1.797.18'
$ws.Range("E13").Value = '  -1.15%  '

$ws.Range("D14").Value = '1.570.28'
$ws.Range("E14").Value = '  -1.24%  '

$ws.Range("D16").Value = '28.536.61'
$ws.Range("E16").Value = '  +0.95%  '

$ws.Range("E17").Value = '  -1.95%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.34'
$ws.Range("E18").Value = '  -1.33%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '230.56'
$ws.Range("E19").Value = '  +1.15%  '

$ws.Range("E20").Value = '  -1.52%  '

$ws.Range("E21").Value = '  -2.58%  '

$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("E23").Value = '  -5.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.11'
$ws.Range("E24").Value = '  -2.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.15'
$ws.Range("E25").Value = '  +10.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.72'
$ws.Range("E26").Value = '  -0.13%  '

$ws.Range("E27").Value = '  -1.40%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.44'
$ws.Range("E28").Value = '  -2.21%  '

$ws.Range("E29").Value = '  -3.21%  '

$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0486'
$ws.Range("E31").Value = '  +2.51%  '

$ws.Range("E32").Value = '  -2.58%  '

$ws.Range("E33").Value = '  -1.11%  '

$ws.Range("E34").Value = '  -1.87%  '

$ws.Range("E35").Value = '  -0.86%  '

$ws.Range("E36").Value = '  +1.37%  '

$ws.Range("E37").Value = '  -3.44%  '

$ws.Range("E38").Value = '  +0.94%  '

$ws.Range("E39").Value = '  +3.32%  '

$ws.Range("E40").Value = '  -0.48%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.525'
$ws.Range("E41").Value = '  -2.89%  '

$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.790'
$ws.Range("E43").Value = '  -2.91%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.89'
$ws.Range("E44").Value = '  +0.47%  '

$ws.Range("E45").Value = '  +2.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.968'
$ws.Range("E47").Value = '  -1.74%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.02'
$ws.Range("E48").Value = '  -2.06%  '

$ws.Range("D49").Value = '1.709.45'
$ws.Range("E49").Value = '  -1.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.54'
$ws.Range("E50").Value = '  -1.28%  '

$ws.Range("E51").Value = '  -0.77%  '
